# Updated cryptos list (price + 1h volume change) as scraped from coinranking.com.
# Values are assigned with a leading apostrophe so Excel stores them as literal
# text (preserving thousands-dot price notation like "68.747.37" and the
# padded percentage strings like "  -0.69%  ") instead of auto-converting them
# to numbers, while keeping the cells' original "General" number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.747.37"
$ws.Range("E2").Value = "'  -0.69%  "
$ws.Range("D3").Value = "'3.836.22"
$ws.Range("E3").Value = "'  +2.22%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'600.90"
$ws.Range("E5").Value = "'  -0.24%  "
$ws.Range("D6").Value = "'161.38"
$ws.Range("E6").Value = "'  -3.43%  "
$ws.Range("D7").Value = "'3.833.66"
$ws.Range("E7").Value = "'  +2.21%  "
$ws.Range("E8").Value = "'  +0.00%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "'  -1.60%  "
$ws.Range("E10").Value = "'  -1.17%  "
$ws.Range("E11").Value = "'  -1.52%  "
$ws.Range("E12").Value = "'  -0.55%  "
$ws.Range("D13").Value = "'36.80"
$ws.Range("E13").Value = "'  -3.38%  "
$ws.Range("E14").Value = "'  -2.41%  "
$ws.Range("D15").Value = "'4.482.76"
$ws.Range("E15").Value = "'  +2.30%  "
$ws.Range("D16").Value = "'3.829.17"
$ws.Range("E16").Value = "'  +2.06%  "
$ws.Range("D17").Value = "'68.894.08"
$ws.Range("E17").Value = "'  -0.47%  "
$ws.Range("D18").Value = "'7.50"
$ws.Range("E18").Value = "'  +1.49%  "
$ws.Range("E19").Value = "'  -0.20%  "
$ws.Range("D20").Value = "'11.34"
$ws.Range("E20").Value = "'  +1.38%  "
$ws.Range("D21").Value = "'17.10"
$ws.Range("E21").Value = "'  -1.85%  "
$ws.Range("D22").Value = "'483.59"
$ws.Range("E22").Value = "'  -2.29%  "
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "'  -1.53%  "
$ws.Range("D24").Value = "'0.0000158"
$ws.Range("E24").Value = "'  +3.16%  "
$ws.Range("D25").Value = "'83.91"
$ws.Range("E25").Value = "'  -1.18%  "
$ws.Range("E26").Value = "'  -2.92%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "'  -2.07%  "
$ws.Range("E28").Value = "'  -0.18%  "
$ws.Range("D29").Value = "'9.97"
$ws.Range("E29").Value = "'  -1.44%  "
$ws.Range("D30").Value = "'2.95"
$ws.Range("E30").Value = "'  -1.36%  "
$ws.Range("E31").Value = "'  -2.09%  "
$ws.Range("D32").Value = "'3.988.55"
$ws.Range("E32").Value = "'  +2.35%  "
$ws.Range("D33").Value = "'2.37"
$ws.Range("E33").Value = "'  -4.11%  "
$ws.Range("D34").Value = "'32.08"
$ws.Range("E34").Value = "'  +1.36%  "
$ws.Range("D35").Value = "'3.784.36"
$ws.Range("E35").Value = "'  +2.52%  "
$ws.Range("D36").Value = "'0.106"
$ws.Range("E36").Value = "'  -2.01%  "
$ws.Range("E37").Value = "'  +0.90%  "
$ws.Range("D38").Value = "'0.140"
$ws.Range("E38").Value = "'  +3.29%  "
$ws.Range("E39").Value = "'  -1.66%  "
$ws.Range("E40").Value = "'  +0.00%  "
$ws.Range("E41").Value = "'  -1.97%  "
$ws.Range("D42").Value = "'436.96"
$ws.Range("E42").Value = "'  +1.43%  "
$ws.Range("E43").Value = "'  -1.13%  "
$ws.Range("D44").Value = "'48.49"
$ws.Range("E44").Value = "'  -0.69%  "
$ws.Range("D45").Value = "'1.98"
$ws.Range("E45").Value = "'  -0.85%  "
$ws.Range("D47").Value = "'8.38"
$ws.Range("E47").Value = "'  -1.44%  "
$ws.Range("B48").Value = "'Monero"
$ws.Range("C48").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'142.99"
$ws.Range("E48").Value = "'  +1.23%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'26.17"
$ws.Range("E49").Value = "'  +12.02%  "
$ws.Range("D50").Value = "'2.821.14"
$ws.Range("E50").Value = "'  +1.00%  "
$ws.Range("E51").Value = "'  +1.86%  "
